# Electricite de Strasbourg - Tarif Bleu / Option Base
# Update of 2026-02-01 prices:
#  - PART_FIXE_HT header renamed to PART_FIXE_TTC_MOIS
#  - Whole column F (PART_VARIABLE_HT, never populated) removed; old column G
#    (PART_VARIABLE_TTC data) shifts left into column F
#  - New pricing rows 93-101 (01/02/2026 onward) get their D/E/F values filled in

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column D header string (PART_FIXE_HT -> PART_FIXE_TTC_MOIS)
$ws.Range("D1").Value = "PART_FIXE_TTC_MOIS"

# Remove the (empty) PART_VARIABLE_HT column; everything to its right shifts left
$ws.Columns("F").Delete()

# Fill in the new tariff data for 01/02/2026 (rows 93-101)
$ws.Range("D93").Value = 12.07
$ws.Range("E93").Formula = '=IF(D93="","",12*D93)'
$ws.Range("F93").Value = 0.19398

$ws.Range("D94").Value = 15.74
$ws.Range("E94").Formula = '=IF(D94="","",12*D94)'
$ws.Range("F94").Value = 0.19398

$ws.Range("D95").Value = 19.69
$ws.Range("E95").Formula = '=IF(D95="","",12*D95)'
$ws.Range("F95").Value = 0.19266

$ws.Range("D96").Value = 23.49
$ws.Range("E96").Formula = '=IF(D96="","",12*D96)'
$ws.Range("F96").Value = 0.19266

$ws.Range("D97").Value = 27.06
$ws.Range("E97").Formula = '=IF(D97="","",12*D97)'
$ws.Range("F97").Value = 0.19266

$ws.Range("D98").Value = 30.75
$ws.Range("E98").Formula = '=IF(D98="","",12*D98)'
$ws.Range("F98").Value = 0.19266

$ws.Range("D99").Value = 38.59
$ws.Range("E99").Formula = '=IF(D99="","",12*D99)'
$ws.Range("F99").Value = 0.19266

$ws.Range("D100").Value = 45.8
$ws.Range("E100").Formula = '=IF(D100="","",12*D100)'
$ws.Range("F100").Value = 0.19266

$ws.Range("D101").Value = 53.06
$ws.Range("E101").Formula = '=IF(D101="","",12*D101)'
$ws.Range("F101").Value = 0.19266

# Update the frozen-pane / selection bookmarks to match the new scroll position
$ws.Application.ActiveWindow.ScrollRow = 77
$ws.Range("F94").Select()
